$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Block 1: "Clock Gating All Design" (rows 18-20)
# Copy the formatting of the last existing block (rows 14-16, "Clock Gating
# Inputs Design") into rows 18-20, then overwrite the text/values that differ.
# ---------------------------------------------------------------------------
$ws.Range("B14:F16").Copy($ws.Range("B18"))

$ws.Range("B18").Value = "Clock Gating All Design"

$ws.Range("B20").Value = 0.0014888002770021599
$ws.Range("C20").Value = 0.00801691412925720041
$ws.Range("D20").Value = 0.00347015284933149988
$ws.Range("F20").Formula = "=B20+C20+D20"

# ---------------------------------------------------------------------------
# Block 2: "Hybrid Clock Gating All and Registering Design" (rows 22-24)
# Same approach: duplicate formatting from the block above, then clear the
# data cells (this new block ships with no figures yet) and fix the formula.
# ---------------------------------------------------------------------------
$ws.Range("B18:F20").Copy($ws.Range("B22"))

$ws.Range("B22").Value = "Hybrid Clock Gating All and Registering Design"

$ws.Range("B24").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("D24").ClearContents()
$ws.Range("F24").Formula = "=B24+C24+D24"

# ---------------------------------------------------------------------------
# Column widths: columns D:F should all end up at width 12 (D was narrower).
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 11.17

# ---------------------------------------------------------------------------
# Selection moves to R11 in the saved view.
# ---------------------------------------------------------------------------
$ws.Range("R11").Select()
